$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 30.319878
$ws.Range("H2").Value = 90.95963399999999
$ws.Range("I2").Value = 0.0886902518702035
$ws.Range("J2").Value = 0.08869025187020349
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 3000.052600179672
$ws.Range("R2").Value = 27000.47340161704
$ws.Range("S2").Value = 0.0186073008097411
$ws.Range("T2").Value = 0.0186073008097411

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 30.319878
$ws.Range("H3").Value = 90.95963399999999
$ws.Range("I3").Value = 0.0886902518702035
$ws.Range("J3").Value = 0.08869025187020349
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 4942.329178652582
$ws.Range("R3").Value = 44480.96260787323
$ws.Range("S3").Value = 0.03065393110855511
$ws.Range("T3").Value = 0.03065393110855511

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 30.319878
$ws.Range("H4").Value = 90.95963399999999
$ws.Range("I4").Value = 0.0886902518702035
$ws.Range("J4").Value = 0.08869025187020349
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 1982.80197580832
$ws.Range("R4").Value = 17845.21778227488
$ws.Range("S4").Value = 0.01229798197798426
$ws.Range("T4").Value = 0.01229798197798426

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 30.319878
$ws.Range("H5").Value = 90.95963399999999
$ws.Range("I5").Value = 0.0886902518702035
$ws.Range("J5").Value = 0.08869025187020349
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 4374.33359365214
$ws.Range("R5").Value = 39369.00234286926
$ws.Range("S5").Value = 0.02713103797392302
$ws.Range("T5").Value = 0.02713103797392302

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 107.5357766666667
$ws.Range("H6").Value = 322.60733
$ws.Range("I6").Value = 0.314558492538282
$ws.Range("J6").Value = 0.3145584925382819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 10640.31281396231
$ws.Range("R6").Value = 95762.81532566076
$ws.Range("S6").Value = 0.06599467663576369
$ws.Range("T6").Value = 0.06599467663576368

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 107.5357766666667
$ws.Range("H7").Value = 322.60733
$ws.Range("I7").Value = 0.314558492538282
$ws.Range("J7").Value = 0.3145584925382819
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 17529.0021539247
$ws.Range("R7").Value = 157761.0193853223
$ws.Range("S7").Value = 0.1087205657504614
$ws.Range("T7").Value = 0.1087205657504614

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 107.5357766666667
$ws.Range("H8").Value = 322.60733
$ws.Range("I8").Value = 0.314558492538282
$ws.Range("J8").Value = 0.3145584925382819
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 7032.421121376179
$ws.Range("R8").Value = 63291.79009238561
$ws.Range("S8").Value = 0.04361736031507802
$ws.Range("T8").Value = 0.04361736031507802

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 107.5357766666667
$ws.Range("H9").Value = 322.60733
$ws.Range("I9").Value = 0.314558492538282
$ws.Range("J9").Value = 0.3145584925382819
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 15514.48724142208
$ws.Range("R9").Value = 139630.3851727987
$ws.Range("S9").Value = 0.09622588983697886
$ws.Range("T9").Value = 0.09622588983697886

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 189.9662756666667
$ws.Range("H10").Value = 569.898827
$ws.Range("I10").Value = 0.5556802318175943
$ws.Range("J10").Value = 0.5556802318175943
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 18796.54064769758
$ws.Range("R10").Value = 169168.8658292782
$ws.Range("S10").Value = 0.1165822512556241
$ws.Range("T10").Value = 0.1165822512556241

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 189.9662756666667
$ws.Range("H11").Value = 569.898827
$ws.Range("I11").Value = 0.5556802318175943
$ws.Range("J11").Value = 0.5556802318175943
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 30965.68750003963
$ws.Range("R11").Value = 278691.1875003567
$ws.Range("S11").Value = 0.192059253247483
$ws.Range("T11").Value = 0.192059253247483

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 189.9662756666667
$ws.Range("H12").Value = 569.898827
$ws.Range("I12").Value = 0.5556802318175943
$ws.Range("J12").Value = 0.5556802318175943
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 12423.05482656674
$ws.Range("R12").Value = 111807.4934391007
$ws.Range("S12").Value = 0.07705182173138878
$ws.Range("T12").Value = 0.07705182173138879

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 189.9662756666667
$ws.Range("H13").Value = 569.898827
$ws.Range("I13").Value = 0.5556802318175943
$ws.Range("J13").Value = 0.5556802318175943
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 27406.96586278095
$ws.Range("R13").Value = 246662.6927650285
$ws.Range("S13").Value = 0.1699869055830984
$ws.Range("T13").Value = 0.1699869055830984

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.040646
$ws.Range("H14").Value = 42.121938
$ws.Range("I14").Value = 0.04107102377392038
$ws.Range("J14").Value = 0.04107102377392038
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 1389.275924543704
$ws.Range("R14").Value = 12503.48332089334
$ws.Range("S14").Value = 0.008616740597870751
$ws.Range("T14").Value = 0.008616740597870752

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.040646
$ws.Range("H15").Value = 42.121938
$ws.Range("I15").Value = 0.04107102377392038
$ws.Range("J15").Value = 0.04107102377392038
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 2288.712850788241
$ws.Range("R15").Value = 20598.41565709417
$ws.Range("S15").Value = 0.01419534060142359
$ws.Range("T15").Value = 0.01419534060142359

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.040646
$ws.Range("H16").Value = 42.121938
$ws.Range("I16").Value = 0.04107102377392038
$ws.Range("J16").Value = 0.04107102377392038
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 918.203583484907
$ws.Range("R16").Value = 8263.832251364161
$ws.Range("S16").Value = 0.00569499690820843
$ws.Range("T16").Value = 0.005694996908208431

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.040646
$ws.Range("H17").Value = 42.121938
$ws.Range("I17").Value = 0.04107102377392038
$ws.Range("J17").Value = 0.04107102377392038
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 2025.683265426647
$ws.Range("R17").Value = 18231.14938883982
$ws.Range("S17").Value = 0.0125639456664176
$ws.Range("T17").Value = 0.01256394566641761
